$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the first empty row after current data
$lastRow = $ws.Cells(1048576, 1).End(-4162).Row  # xlUp = -4162
$newRow = $lastRow + 1

$newData = @(
    @(114, "customer_targeting", "personal", "project", "https://github.com/likarajo/customer_targeting"),
    @(115, "blackfriday_shopping", "personal", "project", "https://github.com/likarajo/blackfriday_shopping"),
    @(116, "text_summarization", "personal", "project", "https://github.com/likarajo/text_summarization"),
    @(117, "wine_quality", "personal", "project", "https://github.com/likarajo/wine_quality"),
    @(118, "shopping_trends", "personal", "project", "https://github.com/likarajo/shopping_trends"),
    @(119, "spam_sms", "personal", "project", "https://github.com/likarajo/spam_sms"),
    @(120, "glass_type", "personal", "project", "https://github.com/likarajo/glass_type"),
    @(121, "currencynote_authenticity", "personal", "project", "https://github.com/likarajo/currencynote_authenticity"),
    @(122, "weather_predcition", "personal", "project", "https://github.com/likarajo/weather_prediction")
)

$r = $newRow
foreach ($row in $newData) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $r = $r + 1
}

$lastRowFinal = $r - 1

# Re-sort the whole data range A2:G(lastRowFinal) by column B ascending
$sortRange = $ws.Range("A2:G$lastRowFinal")
$keyRange = $ws.Range("B2:B$lastRowFinal")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4142  # xlNo
$ws.Sort.Apply()

# Fix up the A column (Sno) to be sequential 1..N
for ($i = 2; $i -le $lastRowFinal; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

$ws.Range("G15").Select()
